$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.644.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.741.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.87%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '116.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '332.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.31%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.572'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.02%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0835'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.43%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.03'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('E13').Value = '  +2.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.55'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.175.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.750.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.879'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.570.48'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.11'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.80'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.28%  '
$ws.Range('E22').Value = '  +2.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '279.45'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.65'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.68'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.85%  '
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('E30').Value = '  +2.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.86'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('E32').Value = '  +2.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0817'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.94'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.52%  '
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.13'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '128.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0345'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.21%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.46%  '
$ws.Range('E43').Value = '  +2.64%  '
$ws.Range('E44').Value = '  +5.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.097.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.14%  '
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.51'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.36%  '
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.85'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.67%  '
